$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at E so the existing "fantasy points" column (E) shifts to G
$ws.Range("E1:F1").EntireColumn.Insert()

# New headers for the inserted columns
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"

# Match the header formatting (bold/border/centered) used by the other header cells
$ws.Range("G1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)

# Fill in the new data for rows 2-12
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.416666666666667
    $ws.Cells.Item($r, 6).Value = 244
}

Write-Output "done"
